$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 2 corresponds to file
# "11d03539-d425-4e32-b99a-31afb4d274be...zh-cn.xlf" - refresh its handoff/
# handback datetimes for a new report generation cycle.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-30 10:09:42"
$wsZh.Range("H2").Value = "2016-03-30 10:10:29"

# de-de sheet: row 2 corresponds to file
# "11d03539-d425-4e32-b99a-31afb4d274be...de-de.xlf" - same refresh.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-30 10:09:52"
$wsDe.Range("H2").Value = "2016-03-30 10:10:46"
